# Edit script: apply the "pasta rotini bowl" addition to the workbook.
# Sheet1 = researchMeasures, Sheet2 = dataDictionary, Sheet3 = NutritionalData
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)
# --- researchMeasures row 46: extend the day's note + food log with the pasta bowl ---
$note46 = 'Woke up at 520 am, and felt tired. Had a couple cups of coffee, didn''t clean up pet messes because the roommate did when he got home, and I made the babies their food. For breakfast around 7:30 am I had 2 scrambled eggs with 1 tbs sour cream and 1 avocado and 2 corn tortilla and mozzarella cheese quesadillas. I also decided to squeeze my waist trimmer in to the 2nd string of grommets or eyes to fasten which is approximately an inch smaller than the 31" I have been wearing since mid January. It didn''t give me any problems at work, and might be another reason not too hungry, but I also ate a big breakfast this morning. For lunch I had a banana and 3 mandarin oranges and my 3rd cup of coffee from the work keurig instead of my instant coffee I didn''t pack any carbs for lunch was running late as I did some of my ch2 hw for genetics and on my lunch break while drinking my coffee I went to the Dollar Tree to get some coconut oil and disposable oil bottles and some plates, bowls, and mugs and try out their foot scrub and some battery lit candles. After work had a banana when I got home and made a cup of my 4th cup of coffee with 1/4 tbs instant nescafe cofffe, 1 tbs cocoa powder 1 tbs honey and 1/2 tbs coconut oil. This version tastes better than without coconut oil. I plan on cooking up the beyond beef burger patties but with peppers and the zuchini in the fridge and one of the pastas. Still not sure. Going to sip my cocoa while I work on the rest of the ch2 hw for genetics, then do some studying for both courses. Need to get the rent money order due tomorrow. But will probably get it tommorrow at my bank or at the grocery store. I actually ended up making the pasta then having a bowl of rotini 1 pkg with 1 pkg of beyond meat for 2 burgers, 1 yellow bell pepper, 1 zucchini, and Prego 3 cheese pasta sauce. I had a whole avocado with the pasta and shared some with the babies. No cheese on this bowl. Then started my genetics ch2 hw. The cocoa with coffee sort of gave me the onset of a headache at first so I made the pasta and ate before doing the homework. After completing the genetics week 2 ch2 homework, I reviewed the power point slides for chemistry 1A week 2 and then decided to watch a movie, the first movie trailer looked good, so picked ''Greenland,'' a movie about the end of the world but only the best and brightest professions like structural engineers and doctors with their immediate family have to fight to take the plane to the bunkers in Greenland for saving the best of humanity while billions of the others die from a large asteroid and its many parts that kill much of the plant and animal life on Earth. Pretty good movie. The little boy''s constant sad, crybaby face of despair was annoying. Because it was like the only face he had on the whole time except when actually saying something about death flashing before your eyes. When you watch other movies with little helpless characters they are usually lively but his emotion the whole time was of a helpless and useless little diabetic kid and I know its mean to say. Its sad really, but he overkills it. Maybe a couple shots of him with that face, but the running around and time constraints weren''t really translated to me to make me feel it like in other end of the world movies. And the way the 7 year old boy stayed in the car when he could have unbuckled himself and got out when his mom was forced out so the strangers could have their bracelets is really where the sad and despaired and pathetic face got annoying. He could have at least got out of the car. But the drama had to be a string of events that are obstacles for the family getting saved. Going to do more review tomorrow. Work out the chapter 2 chemistry worksheet with the answers on the slide first thing in the morning and look at the genetics study sheet as well for chapter 1 and see the items to do for week 3. Bed time around 10:15 pm.'
$food46 = '2 eggs
(140	10	3	12	0	0	140)
1 tbs olive oil
(120	14	2	0	0	0	0)
1 tbs sourcream
(60	5	3.5	1	2	0	15)
1 avocado
(322	29	4	4	17	18	14)
2 corn tortilla quesadillas
4 corn tortillas Guerrero brand
(200	2	0	4	42	4	40)
1/3 cup mozzarella cheese
(106.7	 6.7	4.7	8	1.3	0	253.3)
3 mandarin oranges
(120	0.6	0.3	1.8	30.3	4.2	6)
2 bananas
(210	0	0	2	54	6	2)
1 tbs cocoa powder
(10	0.5	0	1	3	1	0)
1/2 tbs coconut oil
(60	7	6.5	0	0	0	0)
1 tbs honey
(60	0	0	0	17	0	0)

bowl pasta rotini with beyond burger meat prego 3 cheese 
(663.00	25.70	5.95	36.87	76.10	13.10	1035.80)
1 avocado
(322	29	4	4	17	18	14)
4 corn tortillas Guerrero brand
(200	2	0	4	42	4	40)
1/3 cup mozzarella cheese
(106.7	 6.7	4.7	8	1.3	0	253.3)
=140+120+60+322+200+106.7+120+210+10+60+60+663+322+200+106.7
=10+14+5+29+2+6.7+0.6+0+0l5+7+0+25.7+29+2+6.7
=3+2+3.5+4+0+4.7+0.3+0+0+6.5+0+5.95+4+0+4.7
=12+0+1+4+4+8+1.8+2+1+0+0+36.87+4+4+8
=0+0+2+17+42+1.3+30.3+54+3+0+17+76.10+17+42+1.3
=0+0+0+18+4+0+4.2+6+1+0+0+13.10+18+4+0
=140+0+15+14+40+253.3+6+2+0+0+0+1035.8+14+40+253.3
'
$ws1.Range("Z46").Value2 = $note46
$ws1.Range("AA46").Value2 = $food46
$ws1.Range("AB46").Formula = "=140+120+60+322+200+106.7+120+210+10+60+60+663+322+200+106.7"
$ws1.Range("AC46").Formula = "=10+14+5+29+2+6.7+0.6+0+0.5+7+0+25.7+29+2+6.7"
$ws1.Range("AD46").Formula = "=3+2+3.5+4+0+4.7+0.3+0+0+6.5+0+5.95+4+0+4.7"
$ws1.Range("AE46").Formula = "=12+0+1+4+4+8+1.8+2+1+0+0+36.87+4+4+8"
$ws1.Range("AF46").Formula = "=0+0+2+17+42+1.3+30.3+54+3+0+17+76.1+17+42+1.3"
$ws1.Range("AG46").Formula = "=0+0+0+18+4+0+4.2+6+1+0+0+13.1+18+4+0"
$ws1.Range("AH46").Formula = "=140+0+15+14+40+253.3+6+2+0+0+0+1035.8+14+40+253.3"

# --- researchMeasures sheet view: scroll to column U, select AA46 ---
$ws1.Application.ActiveWindow.ScrollColumn = 21
$ws1.Range("AA46").Select()

# --- NutritionalData: add the new pasta-rotini summary rows 120-122 ---
$ws3.Range("A120").Value2 = 'Rotini red lentil barilla brand 1 pkg is 4.5 servings of 2 oz or 3.5 oz in 2.5 servings'
$ws3.Range("B120").Formula = "=330*2.5"
$ws3.Range("C120").Formula = "=2.5+2.5"
$ws3.Range("D120").Formula = "=0.5*2.5"
$ws3.Range("E120").Formula = "=23*2.5"
$ws3.Range("F120").Formula = "=61*2.5"
$ws3.Range("G120").Formula = "=11*2.5"
$ws3.Range("H120").Formula = "=0"

$ws3.Range("A121").Value2 = 'pasta rotini with 3 cheese Prego/beyond meat burgers (2)/1 yellow pepper/1 zuchini/2 tbs olive oil, makes about 3 bowls of pasta'
$ws3.Range("B121").Formula = "=SUM(B120,B119,B114*5,B49*2,B40,B116)"
$ws3.Range("C121").Formula = "=SUM(C120,C119,C114*5,C49*2,C40,C116)"
$ws3.Range("D121").Formula = "=SUM(D120,D119,D114*5,D49*2,D40,D116)"
$ws3.Range("E121").Formula = "=SUM(E120,E119,E114*5,E49*2,E40,E116)"
$ws3.Range("F121").Formula = "=SUM(F120,F119,F114*5,F49*2,F40,F116)"
$ws3.Range("G121").Formula = "=SUM(G120,G119,G114*5,G49*2,G40,G116)"
$ws3.Range("H121").Formula = "=SUM(H120,H119,H114*5,H49*2,H40,H116)"

$ws3.Range("A122").Value2 = '1 bowl pasta rotini with beyond burger meat'
$ws3.Range("B122").Formula = "=B121/3"
$ws3.Range("C122").Formula = "=C121/3"
$ws3.Range("D122").Formula = "=D121/3"
$ws3.Range("E122").Formula = "=E121/3"
$ws3.Range("F122").Formula = "=F121/3"
$ws3.Range("G122").Formula = "=G121/3"
$ws3.Range("H122").Formula = "=H121/3"

# --- match formatting used by the rest of the nutrition-summary block ---
$ws3.Range("A120:A122").HorizontalAlignment = -4131
$ws3.Range("A120:A122").VerticalAlignment = -4160
$ws3.Range("B122:H122").NumberFormat = "0.00"

# --- NutritionalData sheet view: scroll down, select the new averages row ---
$ws3.Range("B122:H122").Select()
